# Added Upload Profile Test with Excel connect
$wb = $excel.ActiveWorkbook

# Rename Sheet2 to reflect the new upload-profile test
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "test_invalid_profile_upload"

# Populate the header row
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "Upload File Detail"
$ws2.Range("D1").Value = "Expected Error"

# Fill in the user credentials columns first
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "admin123"
$ws2.Range("A3").Value = "Admin"
$ws2.Range("B3").Value = "admin123"

# Then the upload-file-detail column
$ws2.Range("C2").Value = "C:\Mine\Balaji-Profile_2023.pdf"
$ws2.Range("C3").Value = "C:\Mine\iFuture.txt"

# Then the expected-error column (same message both rows)
$ws2.Range("D2").Value = "File type not allowed"
$ws2.Range("D3").Value = "File type not allowed"

# Match column widths to the target layout
$ws2.Columns.Item(2).ColumnWidth = 9
$ws2.Columns.Item(3).ColumnWidth = 27.44140625
$ws2.Columns.Item(4).ColumnWidth = 18.109375

# Page setup (portrait) for the new sheet
$ws2.PageSetup.Orientation = 1

# Selection on the new sheet
$ws2.Range("D2:D3").Select()

# Make this sheet the active one (sets tabSelected/activeTab)
$ws2.Activate()
